$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.07
$ws.Range("O2").Value = 1.41
$ws.Range("P2").Value = 2.62

$ws.Range("L4").Value = 7
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 8.5
$ws.Range("Z4").Value = 10
$ws.Range("AF4").Value = 81
$ws.Range("AG4").Value = 12
$ws.Range("AL4").Value = 51
$ws.Range("AZ4").Value = 151
$ws.Range("BA4").Value = 201

$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25

$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 13

$ws.Range("G10").Value = 2.45
$ws.Range("M10").Value = 1.13
$ws.Range("N10").Value = 6
$ws.Range("U10").Value = 2.2
$ws.Range("V10").Value = 1.62
$ws.Range("Y10").Value = 11
$ws.Range("AC10").Value = 6
$ws.Range("AH10").Value = 13

$ws.Range("G14").Value = 5.25
$ws.Range("H14").Value = 3.6
$ws.Range("I14").Value = 1.55
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 11
$ws.Range("W14").Value = 13
$ws.Range("Z14").Value = 51
$ws.Range("AC14").Value = 10
$ws.Range("AF14").Value = 51
$ws.Range("AN14").Value = 7
$ws.Range("AO14").Value = 29
$ws.Range("AW14").Value = 3.5

$ws.Range("M16").Value = 1.03
$ws.Range("O16").Value = 1.25

$ws.Range("G17").Value = 1.85
$ws.Range("I17").Value = 3.9
$ws.Range("J17").Value = 2.6
$ws.Range("N17").Value = 10
$ws.Range("Q17").Value = 2.05
$ws.Range("R17").Value = 1.8
$ws.Range("S17").Value = 1.44
$ws.Range("T17").Value = 2.63
$ws.Range("X17").Value = 9
$ws.Range("Y17").Value = 9
$ws.Range("AA17").Value = 17
$ws.Range("AG17").Value = 11
$ws.Range("AI17").Value = 13
$ws.Range("AM17").Value = 301
$ws.Range("AT17").Value = 2.63
$ws.Range("AX17").Value = 21
